$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-28 Sunday" "2024-07-29 Monday"

Replace-Text "229÷8=28, 5" "691÷7=98, 5"
Replace-Text "230÷2=115, 0" "825÷5=165, 0"
Replace-Text "157÷7=22, 3" "964÷9=107, 1"
Replace-Text "744÷8=93, 0" "211÷7=30, 1"
Replace-Text "817÷8=102, 1" "971÷5=194, 1"

Replace-Text "771÷9=85, 6" "515÷8=64, 3"
Replace-Text "466÷9=51, 7" "671÷8=83, 7"
Replace-Text "872÷4=218, 0" "918÷3=306, 0"
Replace-Text "373÷3=124, 1" "165÷7=23, 4"
Replace-Text "493÷7=70, 3" "920÷9=102, 2"

Replace-Text "952÷8=119, 0" "522÷6=87, 0"
Replace-Text "783÷2=391, 1" "140÷8=17, 4"
Replace-Text "781÷3=260, 1" "409÷7=58, 3"
Replace-Text "580÷5=116, 0" "686÷4=171, 2"
Replace-Text "980÷2=490, 0" "518÷2=259, 0"

Replace-Text "952÷6=158, 4" "390÷4=97, 2"
Replace-Text "333÷5=66, 3" "939÷3=313, 0"
Replace-Text "635÷4=158, 3" "515÷6=85, 5"
Replace-Text "436÷8=54, 4" "311÷2=155, 1"
Replace-Text "450÷3=150, 0" "382÷4=95, 2"

Replace-Text "641÷9=71, 2" "866÷4=216, 2"
Replace-Text "639÷8=79, 7" "966÷9=107, 3"
Replace-Text "655÷5=131, 0" "585÷3=195, 0"
Replace-Text "545÷3=181, 2" "886÷7=126, 4"
Replace-Text "185÷3=61, 2" "298÷4=74, 2"
